$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C (shifts old C..I to D..J), preserving styles/formatting
$ws.Columns.Item(3).Insert()

# Header for new Industry column
$ws.Cells.Item(1,3).Value = "Industry"

# Populate Industry values for each data row
$ws.Cells.Item(2,3).Value = "Banks"
$ws.Cells.Item(3,3).Value = "IT - Software"
$ws.Cells.Item(4,3).Value = "Banks"
$ws.Cells.Item(5,3).Value = "Finance"
$ws.Cells.Item(6,3).Value = "Telecom - Services"
$ws.Cells.Item(7,3).Value = "IT - Software"
$ws.Cells.Item(8,3).Value = "Healthcare Services"
$ws.Cells.Item(9,3).Value = "Banks"
$ws.Cells.Item(10,3).Value = "IT - Software"
$ws.Cells.Item(11,3).Value = "Consumer Durables"
$ws.Cells.Item(12,3).Value = "Food Products"
$ws.Cells.Item(13,3).Value = "Auto Components"
$ws.Cells.Item(14,3).Value = "Retailing"
$ws.Cells.Item(15,3).Value = "Food Products"
$ws.Cells.Item(16,3).Value = "Automobiles"
$ws.Cells.Item(17,3).Value = "Finance"
$ws.Cells.Item(18,3).Value = "Finance"
$ws.Cells.Item(19,3).Value = "IT - Software"
$ws.Cells.Item(20,3).Value = "Pharmaceuticals & Biotechnology"
$ws.Cells.Item(21,3).Value = "Healthcare Equipment & Supplies"
$ws.Cells.Item(22,3).Value = "Chemicals & Petrochemicals"
$ws.Cells.Item(23,3).Value = "Realty"
$ws.Cells.Item(24,3).Value = "Pharmaceuticals & Biotechnology"
$ws.Cells.Item(25,3).Value = "Realty"
$ws.Cells.Item(26,3).Value = "Industrial Manufacturing"
$ws.Cells.Item(27,3).Value = "Electrical Equipment"
$ws.Cells.Item(28,3).Value = "Automobiles"
$ws.Cells.Item(29,3).Value = "Capital Markets"
$ws.Cells.Item(30,3).Value = "Leisure Services"
$ws.Cells.Item(31,3).Value = "Auto Components"
$ws.Cells.Item(32,3).Value = "Capital Markets"
$ws.Cells.Item(33,3).Value = "Agricultural Food & other Products"
$ws.Cells.Item(34,3).Value = "Banks"
$ws.Cells.Item(35,3).Value = "Financial Technology (Fintech)"
$ws.Cells.Item(36,3).Value = "Healthcare Services"
$ws.Cells.Item(37,3).Value = "Insurance"
$ws.Cells.Item(38,3).Value = "Power"
$ws.Cells.Item(39,3).Value = "Banks"
$ws.Cells.Item(40,3).Value = "Pharmaceuticals & Biotechnology"
$ws.Cells.Item(41,3).Value = "Food Products"
$ws.Cells.Item(42,3).Value = "Food Products"
$ws.Cells.Item(43,3).Value = "Diversified FMCG"
$ws.Cells.Item(44,3).Value = "Banks"
$ws.Cells.Item(45,3).Value = "Electrical Equipment"
$ws.Cells.Item(46,3).Value = "Auto Components"

$ws.Range("A1").Select()
